# 22nd commit of Docker Project - add "Redundancies" slide, tweak "To-do" slide.
$p = $ppt.ActivePresentation

# ------------------------------------------------------------------
# 1. Insert the new "Redundancies" slide at position 6 (right after the
#    "Process of making modules" slide, before "Problems"). Layout 2 =
#    ppLayoutText = Title + Content, same layout used by the sibling
#    slides in this deck.
# ------------------------------------------------------------------
$newSlide = $p.Slides.Add(6, 2)

$title = $newSlide.Shapes.Item(1)
$title.TextFrame.TextRange.Text = "Redundancies"

$body = $newSlide.Shapes.Item(2)
$bodyTr = $body.TextFrame.TextRange
$bodyTr.Text = "The use of Volumes " + [char]8211 + " can be backed up elsewhere " + [char]8211 + " NAS, Offsite, etc."
$bodyTr.InsertAfter("`rThe use of storage containers") | Out-Null
$bodyTr.InsertAfter("`r") | Out-Null
$bodyTr.InsertAfter("`rNo data is theoretically lost should a container ") | Out-Null
$bodyTr.InsertAfter("go down.") | Out-Null

# ------------------------------------------------------------------
# 2. Update the "To-do" slide (now pushed one position later, at the
#    end of the deck): drop the "The use of Volumes" bullet and tweak
#    three remaining lines.
# ------------------------------------------------------------------
$todo = $p.Slides.Item($p.Slides.Count)
$todoBody = $todo.Shapes.Item(2)
$todoTr = $todoBody.TextFrame.TextRange
$todoTr.Text = "Data redundancies via the use of storage containers`rUrbancode (If there is a IBM update) deploy with the use of MySQL`rPuppet Enterprise ?`rThe dreaded Postfix of doom " + [char]8211 + " Properly  "
